# Update numeric values on "Repayment schedule" sheet to reflect the
# revised amortisation schedule (11 testcases on Multi Reschedule).

$wb = $excel.ActiveWorkbook
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Row -> (F, G, H) new values
$rows = @{
    3  = @{ F = 788.49;  G = 9211.51;             H = 100 }
    4  = @{ F = 796.37;  G = 8415.14;             H = 92.12 }
    5  = @{ F = 804.34;  G = 7610.8;              H = 84.15 }
    6  = @{ F = 812.38;  G = 6798.42;             H = 76.11 }
    7  = @{ F = 820.51;  G = 5977.91;             H = 67.98 }
    8  = @{ F = 828.71;  G = 5149.2;              H = 59.78 }
    9  = @{ F = 837;     G = 4312.2;              H = 51.49 }
    10 = @{ F = 845.37;  G = 3466.83;             H = 43.12 }
    11 = @{ F = 853.82;  G = 2613.0100000000002;  H = 34.67 }
    12 = @{ F = 862.36;  G = 1750.65;             H = 26.13 }
    13 = @{ F = 870.98;  G = 879.67;              H = 17.510000000000002 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $wsRepay.Cells.Item($r, 6).Value = $vals.F   # column F - Principal Due
    $wsRepay.Cells.Item($r, 7).Value = $vals.G   # column G - Balance of Loan
    $wsRepay.Cells.Item($r, 8).Value = $vals.H   # column H - Interest
}

# Row 14 has different changes: F, H, K, P (G stays 0)
$wsRepay.Cells.Item(14, 6).Value = 879.67                    # F14
$wsRepay.Cells.Item(14, 8).Value = 8.8000000000000007        # H14
$wsRepay.Cells.Item(14, 11).Value = 888.47                   # K14
$wsRepay.Cells.Item(14, 16).Value = 888.47                   # P14

# --- Selection / active-tab bookkeeping, mirroring the view state left
# behind by whoever ran the testcases (sheet1 -> sheet2 -> sheet3 order,
# ending with "Repayment schedule" as the active tab). ---

$wsInput = $wb.Worksheets.Item("NewLoanInput")
$wsInput.Activate()
$wsInput.Range("A1:B11").Select()

$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEdit.Activate()
$wsEdit.Range("B8").Select()

$wsRepay.Activate()
$wsRepay.Range("K6").Select()
